$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.6191033333333333
$ws.Range("H2").Value = 1.85731
$ws.Range("I2").Value = 0.007929637811697733
$ws.Range("J2").Value = 0.007962490129789305
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 16.67754233333333
$ws.Range("N2").Value = 50.03262699999999
$ws.Range("O2").Value = 0.9535192900707901
$ws.Range("P2").Value = 0.9578676752791928
$ws.Range("Q2").Value = 10.32512205037444
$ws.Range("R2").Value = 92.92609845336999
$ws.Range("S2").Value = 0.007561062616728516
$ws.Range("T2").Value = 0.0076270119100548
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.6191033333333333
$ws.Range("H3").Value = 1.85731
$ws.Range("I3").Value = 0.007929637811697733
$ws.Range("J3").Value = 0.007962490129789305
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.574769
$ws.Range("N3").Value = 1.724307
$ws.Range("O3").Value = 0.0328617561197435
$ws.Range("P3").Value = 0.03301161735036698
$ws.Range("Q3").Value = 0.3558414037966667
$ws.Range("R3").Value = 3.20257263417
$ws.Range("S3").Value = 0.0002605818238859074
$ws.Range("T3").Value = 0.0002628546773206785
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.6191033333333333
$ws.Range("H4").Value = 1.85731
$ws.Range("I4").Value = 0.007929637811697733
$ws.Range("J4").Value = 0.007962490129789305
$ws.Range("K4").Value = 2
$ws.Range("M4").Value = 0.2382025
$ws.Range("N4").Value = 0.476405
$ws.Range("O4").Value = 0.01361895380946642
$ws.Range("P4").Value = 0.009120707370440172
$ws.Range("Q4").Value = 0.1474719617583334
$ws.Range("R4").Value = 0.8848317705500001
$ws.Range("S4").Value = 0.0001079933710833098
$ws.Range("T4").Value = 0.00007262354241382643
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 76.48912033333333
$ws.Range("H5").Value = 229.467361
$ws.Range("I5").Value = 0.9796927073757713
$ws.Range("J5").Value = 0.9837515530909214
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 16.67754233333333
$ws.Range("N5").Value = 50.03262699999999
$ws.Range("O5").Value = 0.9535192900707901
$ws.Range("P5").Value = 0.9578676752791928
$ws.Range("Q5").Value = 1275.650542398594
$ws.Range("R5").Value = 11480.85488158734
$ws.Range("S5").Value = 0.9341558948244757
$ws.Range("T5").Value = 0.9423038132114964
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 76.48912033333333
$ws.Range("H6").Value = 229.467361
$ws.Range("I6").Value = 0.9796927073757713
$ws.Range("J6").Value = 0.9837515530909214
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.574769
$ws.Range("N6").Value = 1.724307
$ws.Range("O6").Value = 0.0328617561197435
$ws.Range("P6").Value = 0.03301161735036698
$ws.Range("Q6").Value = 43.96357520486966
$ws.Range("R6").Value = 395.672176843827
$ws.Range("S6").Value = 0.03219442282207383
$ws.Range("T6").Value = 0.03247522983846672
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 76.48912033333333
$ws.Range("H7").Value = 229.467361
$ws.Range("I7").Value = 0.9796927073757713
$ws.Range("J7").Value = 0.9837515530909214
$ws.Range("K7").Value = 2
$ws.Range("M7").Value = 0.2382025
$ws.Range("N7").Value = 0.476405
$ws.Range("O7").Value = 0.01361895380946642
$ws.Range("P7").Value = 0.009120707370440172
$ws.Range("Q7").Value = 18.21989968620083
$ws.Range("R7").Value = 109.319398117205
$ws.Range("S7").Value = 0.01334238972922173
$ws.Range("T7").Value = 0.008972510040958333
$ws.Range("E8").Value = 2
$ws.Range("G8").Value = 0.9663805
$ws.Range("H8").Value = 1.932761
$ws.Range("I8").Value = 0.01237765481253107
$ws.Range("J8").Value = 0.008285956779289245
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 16.67754233333333
$ws.Range("N8").Value = 50.03262699999999
$ws.Range("O8").Value = 0.9535192900707901
$ws.Range("P8").Value = 0.9578676752791928
$ws.Range("Q8").Value = 16.11685169885783
$ws.Range("R8").Value = 96.70111019314697
$ws.Range("S8").Value = 0.01180233262958592
$ws.Range("T8").Value = 0.007936850157641657
$ws.Range("E9").Value = 2
$ws.Range("G9").Value = 0.9663805
$ws.Range("H9").Value = 1.932761
$ws.Range("I9").Value = 0.01237765481253107
$ws.Range("J9").Value = 0.008285956779289245
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.574769
$ws.Range("N9").Value = 1.724307
$ws.Range("O9").Value = 0.0328617561197435
$ws.Range("P9").Value = 0.03301161735036698
$ws.Range("Q9").Value = 0.5554455536045
$ws.Range("R9").Value = 3.332673321627
$ws.Range("S9").Value = 0.0004067514737837655
$ws.Range("T9").Value = 0.0002735328345795757
$ws.Range("E10").Value = 2
$ws.Range("G10").Value = 0.9663805
$ws.Range("H10").Value = 1.932761
$ws.Range("I10").Value = 0.01237765481253107
$ws.Range("J10").Value = 0.008285956779289245
$ws.Range("K10").Value = 2
$ws.Range("M10").Value = 0.2382025
$ws.Range("N10").Value = 0.476405
$ws.Range("O10").Value = 0.01361895380946642
$ws.Range("P10").Value = 0.009120707370440172
$ws.Range("Q10").Value = 0.23019425105125
$ws.Range("R10").Value = 0.920777004205
$ws.Range("S10").Value = 0.0001685707091613803
$ws.Range("T10").Value = 0.00007557378706801213
